$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Sponsored
Lödha Bhandup | Luxury 2, 2.5 & 3 BHK - Avail Pre Launch Offer
l-bhandup.com
https://www.l-bhandup.com › official
Lödha on LBS Rd an upcoming residential project with the best of amenities & views. Pre Launch Offer. EOI Bookings Open, Pay just...
Brochure & Floor Plans · Price · Connect on Whatsapp · View Project Highlights
Deal: Up to ₹500,000 off Pre Launch Offer'
$ws.Range("B2").Value = 'Rioga Premium Real Estate Advisory LLP'

$ws.Range("A3").Value = 'Sponsored
Lodha Bhandup - New Project Launch in Bhandup
proptigermumbai.com
https://www.proptigermumbai.com
its Big Apartments with 2/3 BHK, its prime Location along the LBS Road, Starts @ ₹ 2.29Cr* The Apartments have 2/3 BHK Configurations & are Designed to Provide Ample Space & Comfort. New Launch Project. Easy Payment Plan. Flexible Payment Plan. Budget Friendly.
Pricing & Floor Plan · Download Brochure · Runwal Group · Platinum Group · Adani Realty'
$ws.Range("B3").Value = 'Locon Solutions Pvt. Ltd.'

$ws.Range("A4").Value = 'Sponsored
New Launch At Bhandup, Mumbai | New Launch At Bhandup
lódháhomz.site
https://www.lódháhomz.site › visit-site › enquire-now
New Launch At Bhandup provides ultra-luxurious residences with the most breathtaking views
E-Brochure · Pricing / Costing · Floor Plan · Project Overview'
$ws.Range("B4").Value = 'DIGITAL RUBIX'

$ws.Range("A5").Value = 'Sponsored
Lodha Prelaunching Bhandup - 2 & 3 BHK Starting ₹2.29 Cr*
prelaunch-projects.in
https://www.prelaunch-projects.in › lodhabhandup › luxuryhomes
Get EOI and Early Bird Benefits, Prime location at LBS Road Bhandup by Lodha. Lodha...'
$ws.Range("B5").Value = 'PRELAUNCH REALTY PRIVATE LIMITED'
$ws.Range("C5").Value = 'India'

$ws.Range("A6").Value = 'Sponsored
Lodha New Launch Bhandup | Luxury 2 & 3 BHK
bhandupnewlaunch.com
https://www.bhandupnewlaunch.com › 2&3bhk › luxury_homes
Pre-book Lodha Bhandup at ₹1.08 Lacs* | Easy Access to Powai & R-City Mall | EOI Open Now! Modern Living at Lodha Bhandup | 10 mins to Eastern Express Hwy | Pre-book at ₹1.08 Lacs* Free Pickup & Drop. Avail Special Offers. Book Now.'
$ws.Range("B6").Value = 'Finwizz Holdings'

$ws.Range("A7").Value = 'Sponsored
Lodha® LBS is Coming Soon
lodhagroup.in
https://www.lodhagroup.in
Lodha® coming soon to LBS — Live an exceptional lifestyle with forest living as Lodha comes to the prime LBS Marg. Lodha LBS along the Mulund-Bhandup-Kanjurmarg corridor ensures seamless connectivity.'
$ws.Range("B7").Value = 'Macrotech Developers Limited'

$ws.Range("A8").Value = 'Sponsored
Lodha Bhandup | Lodha Bhandup New Launch
propertymumbai.co.in
http://www.propertymumbai.co.in › lodha_bhandup
Book 2 & 3 BHK Starts ₹2.29 Cr* & Get EOI & Early Bird Benefits at New Launch Bhandup West. its Big Apartments with 2/3 BHK, its prime Location along the LBS Road, Starts @ ₹ 2.29Cr* Download Brochure.
Price · Floor Plan · Location · Amenities'
$ws.Range("B8").Value = 'NORA GROUP'
